# Update dBu meter calculation.xlsx
#
# 1. Rename the existing "Sheet1" to "11v75" and update its data/formulas
#    for a fixed J5 (total) of 100.
# 2. Duplicate that sheet as "10V" and overwrite it with the 10V-specific
#    measurements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "11v75"

# --- 11v75 sheet -----------------------------------------------------
# J5 becomes a hard-coded total instead of SUM(J6:J15)
$ws.Range("J5").Value = 100

# J7:J15 get new measured values
$ws.Range("J7").Value = 5.49
$ws.Range("J8").Value = 3.83
$ws.Range("J9").Value = 2.74
$ws.Range("J10").Value = 1.91
$ws.Range("J11").Value = 1.37
$ws.Range("J12").Value = 0.98
$ws.Range("J13").Value = 0.68
$ws.Range("J14").Value = 0.47
$ws.Range("J15").Value = 1.18

# K7:K14 change from J.n/$J$5 to SUM(J.n:J15)/J5 (no longer a shared formula)
$ws.Range("K7").Formula = "=SUM(J7:J15)/J5"
$ws.Range("K8").Formula = "=SUM(J8:J15)/J5"
$ws.Range("K9").Formula = "=SUM(J9:J15)/J5"
$ws.Range("K10").Formula = "=SUM(J10:J15)/J5"
$ws.Range("K11").Formula = "=SUM(J11:J15)/J5"
$ws.Range("K12").Formula = "=SUM(J12:J15)/J5"
$ws.Range("K13").Formula = "=SUM(J13:J15)/J5"
$ws.Range("K14").Formula = "=SUM(J14:J15)/J5"
# K15 keeps J15/$J$5 form but re-anchors its shared-formula range
$ws.Range("K15").Formula = "=J15/`$J`$5"

# L7:L15 keep the same formula shape (K.n*$I$5) - re-enter so the shared
# range shrinks from L7:L16 to L7:L15
$ws.Range("L7").Formula = "=K7*`$I`$5"
$ws.Range("L8").Formula = "=K8*`$I`$5"
$ws.Range("L9").Formula = "=K9*`$I`$5"
$ws.Range("L10").Formula = "=K10*`$I`$5"
$ws.Range("L11").Formula = "=K11*`$I`$5"
$ws.Range("L12").Formula = "=K12*`$I`$5"
$ws.Range("L13").Formula = "=K13*`$I`$5"
$ws.Range("L14").Formula = "=K14*`$I`$5"
$ws.Range("L15").Formula = "=K15*`$I`$5"

# J6 becomes a computed remainder (total minus the rest) in bold red text,
# K6/L6 are cleared out (keeping their existing number-format style)
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("J6").Formula = "=J5-SUM(J7:J15)"
$ws.Range("J6").Font.Bold = $true
$ws.Range("J6").Font.Color = 255

# New rows below the table
$ws.Range("J19").Value = 80.6
$ws.Range("J20").Value = 0.75

$ws.Range("I5").Select()

# --- 10V sheet ---------------------------------------------------------
# Duplicate 11v75 (post-edit) so formulas/styles start identical, then
# rename and overwrite with the 10V-specific numbers.
$ws.Copy($null, $ws)
$ws2 = $wb.Worksheets.Item($ws.Index + 1)
$ws2.Name = "10V"

# Clear the extra rows 19-20 that don't belong on this sheet
$ws2.Range("J19").ClearContents()
$ws2.Range("J20").ClearContents()

$ws2.Range("I5").Value = 10
$ws2.Range("J5").Value = 100000

# Stray leftovers that showed up on this sheet from manual editing in the
# source workbook: L6 carries the "Vp" label and I7 picked up the 0.000
# number format even though it keeps the shared I7:I14 formula.
$ws2.Range("L6").Value = "Vp"
$ws2.Range("I7").NumberFormat = "0.000"

$ws2.Range("J7").Value = 6490
$ws2.Range("J8").Value = 4530
$ws2.Range("J9").Value = 3160
$ws2.Range("J10").Value = 2260
$ws2.Range("J11").Value = 1620
$ws2.Range("J12").Value = 1130
$ws2.Range("J13").Value = 787
$ws2.Range("J14").Value = 560
$ws2.Range("J15").Value = 1400

$ws2.Range("J20").Value = 76800
$ws2.Range("J21").Value = 1270
$ws2.Range("J22").Formula = "=J21+J20"

$ws2.Range("E22").Select()
